$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.493.46'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.23%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.807.63'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.12%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.48%  '

# Row 5
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.39%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '306.87'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.50%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4522'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.54%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3600'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.64%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.43'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.50%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07079'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.66%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8899'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.20%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07803'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.37%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.46'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.29%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.775.32'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.45%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.292'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.22%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.322'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.03%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '85.37'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.36%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.006'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.40%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008500'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.79%  '

# Row 20
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.41%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.517.07'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.29%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.20'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.18%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.966'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.20%  '

# Row 24
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.52'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.66%  '

# Row 25
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.000.16'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.92%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.962'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.16%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.36'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.19%  '

# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.33%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.064'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.48%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '112.11'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.70%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.865'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.03%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08694'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.07%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.100'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.92%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.836'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +13.37%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.446'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.22%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7228'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.27%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.104'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.79%  '

# Row 38
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.073'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.01%  '

# Row 39
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01935'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.99%  '

# Row 40
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05104'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.24%  '

# Row 41
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.887'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.68%  '

# Row 42
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5117'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.34%  '

# Row 43
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.784'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.46%  '

# Row 44
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1511'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.64%  '

# Row 45
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.023'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.23%  '

# Row 46
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4667'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.65%  '

# Row 47
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.004'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.39%  '

# Row 48
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.955'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.34%  '

# Row 49
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '100.51'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.53%  '

# Row 50
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.574'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.54%  '

# Row 51
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05981'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.24%  '
